$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

$ws.Range("A28").Value = "Can't Hurt Me"
$ws.Range("B28").Value = "David Goggins"
$ws.Range("E28").Value = "david goggins;no weakness;mental toughness;fitness;strong;navy seals"
$ws.Range("F28").Value = "Audio"

# Copy the date cell's format (style) from the row above, then set the value
$ws.Range("C27").Copy()
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = 43879

$ws.Range("G28").Select()
